$d = $word.ActiveDocument

$d.Content.Find.Execute("47×48=2256", $true, $false, $false, $false, $false, $true, 1, $false, "75×94=7050", 2)
$d.Content.Find.Execute("85×76=6460", $true, $false, $false, $false, $false, $true, 1, $false, "57×98=5586", 2)
$d.Content.Find.Execute("66×46=3036", $true, $false, $false, $false, $false, $true, 1, $false, "47×31=1457", 2)
$d.Content.Find.Execute("23×25=575", $true, $false, $false, $false, $false, $true, 1, $false, "29×93=2697", 2)
$d.Content.Find.Execute("61×79=4819", $true, $false, $false, $false, $false, $true, 1, $false, "77×32=2464", 2)
$d.Content.Find.Execute("86×88=7568", $true, $false, $false, $false, $false, $true, 1, $false, "80×19=1520", 2)
$d.Content.Find.Execute("71×66=4686", $true, $false, $false, $false, $false, $true, 1, $false, "59×26=1534", 2)
$d.Content.Find.Execute("43×42=1806", $true, $false, $false, $false, $false, $true, 1, $false, "89×71=6319", 2)
$d.Content.Find.Execute("95×49=4655", $true, $false, $false, $false, $false, $true, 1, $false, "19×51=969", 2)
$d.Content.Find.Execute("71×96=6816", $true, $false, $false, $false, $false, $true, 1, $false, "82×43=3526", 2)
$d.Content.Find.Execute("92×53=4876", $true, $false, $false, $false, $false, $true, 1, $false, "20×21=420", 2)
$d.Content.Find.Execute("38×71=2698", $true, $false, $false, $false, $false, $true, 1, $false, "16×72=1152", 2)
$d.Content.Find.Execute("69×35=2415", $true, $false, $false, $false, $false, $true, 1, $false, "83×84=6972", 2)
$d.Content.Find.Execute("77×35=2695", $true, $false, $false, $false, $false, $true, 1, $false, "90×76=6840", 2)
$d.Content.Find.Execute("68×93=6324", $true, $false, $false, $false, $false, $true, 1, $false, "72×96=6912", 2)
$d.Content.Find.Execute("60×21=1260", $true, $false, $false, $false, $false, $true, 1, $false, "26×90=2340", 2)
$d.Content.Find.Execute("40×94=3760", $true, $false, $false, $false, $false, $true, 1, $false, "91×31=2821", 2)
$d.Content.Find.Execute("16×57=912", $true, $false, $false, $false, $false, $true, 1, $false, "62×45=2790", 2)
$d.Content.Find.Execute("97×39=3783", $true, $false, $false, $false, $false, $true, 1, $false, "33×37=1221", 2)
$d.Content.Find.Execute("87×55=4785", $true, $false, $false, $false, $false, $true, 1, $false, "11×92=1012", 2)
$d.Content.Find.Execute("28×59=1652", $true, $false, $false, $false, $false, $true, 1, $false, "26×56=1456", 2)
$d.Content.Find.Execute("26×39=1014", $true, $false, $false, $false, $false, $true, 1, $false, "69×42=2898", 2)
$d.Content.Find.Execute("45×89=4005", $true, $false, $false, $false, $false, $true, 1, $false, "19×23=437", 2)
$d.Content.Find.Execute("64×29=1856", $true, $false, $false, $false, $false, $true, 1, $false, "80×56=4480", 2)
$d.Content.Find.Execute("96×87=8352", $true, $false, $false, $false, $false, $true, 1, $false, "11×33=363", 2)
